$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "pants for men sport"
$ws.Range("A2").Value = "mens pouch leggings"
$ws.Range("A3").Value = "spandex leggings boys"
$ws.Range("A4").Value = "rodilleras de basketball"
$ws.Range("A5").Value = "knee sleeves wrestling"
$ws.Range("A6").Value = "knee pads workout"
$ws.Range("A7").Value = "baseball catcher leg guards adult"
$ws.Range("A8").Value = "5 inch seam shorts men"
$ws.Range("A9").Value = "youth softball pants for girls"
$ws.Range("A10").Value = "mens compression workout pants"
$ws.Range("A11").Value = "compression shorts for men"
$ws.Range("A12").Value = "work knee pads for men"
$ws.Range("A13").Value = "padded shorts snowboarding"
$ws.Range("A14").Value = "youth girls yoga pants"
$ws.Range("A15").Value = "youth hockey padded shorts"
$ws.Range("A16").Value = "knee protector sports"
$ws.Range("A17").Value = "kneeling pad construction"
$ws.Range("A18").Value = "knee pad for work"
$ws.Range("A19").Value = "basketballs under"
$ws.Range("A20").Value = "boys baseball pants size 6"
$ws.Range("A21").Value = "sliding short"
$ws.Range("A22").Value = "black knee pads volleyball girls"
$ws.Range("A23").Value = "mens 3/4 pants"
$ws.Range("A24").Value = "football leg sleeves for men"
$ws.Range("A25").Value = "black basketball shorts men"
$ws.Range("A26").Value = "compression football shorts"
$ws.Range("A27").Value = "running leggings men"
$ws.Range("A28").Value = "football leggings boys"
$ws.Range("A29").Value = "baseball leg guards"
$ws.Range("A30").Value = "protective basketball"
$ws.Range("A31").Value = "mountain bike knee pads"
$ws.Range("A32").Value = "weightlifting shorts men"
$ws.Range("A33").Value = "shorts for men basketball"
$ws.Range("A34").Value = "wrestling shorts for boys"
$ws.Range("A35").Value = "baseball items for men"
$ws.Range("A36").Value = "knee pads for biking"
$ws.Range("A37").Value = "eva foam knee pads"
$ws.Range("A38").Value = "mens compression running tights"
$ws.Range("A39").Value = "womens softball pants black"
$ws.Range("A40").Value = "waist guard"
$ws.Range("A41").Value = "bump pads"
$ws.Range("A42").Value = "mens fitness pants"
$ws.Range("A43").Value = "cycling pants for men padded"
$ws.Range("A44").Value = "knee pads work"
$ws.Range("A45").Value = "youth football girdle"
$ws.Range("A46").Value = "bjj knee sleeves"
$ws.Range("A47").Value = "volleyball kneepads black"
$ws.Range("A48").Value = "non slip knee pads"
$ws.Range("A49").Value = "short baseball"
$ws.Range("A50").Value = "knee sleeves for wrestling"
$ws.Range("A51").Value = "knee pads"
$ws.Range("A52").Value = "compression sleeve youth baseball"
$ws.Range("A53").Value = "tights for football"
$ws.Range("A54").Value = "soccer pad"
$ws.Range("A55").Value = "cycling pants for men"
$ws.Range("A56").Value = "dry fit leggings men"
$ws.Range("A57").Value = "compression calf sleeve men basketball"
$ws.Range("A58").Value = "long compression shorts men"
$ws.Range("A59").Value = "compression shorts long men"
$ws.Range("A60").Value = "compression pants and tops for men"
$ws.Range("A61").Value = "leggings knee length"
$ws.Range("A62").Value = "mens softball gear"
$ws.Range("A63").Value = "yoga after knee replacement"
$ws.Range("A64").Value = "wrestling knee sleeve youth"
$ws.Range("A65").Value = "tights compression"
$ws.Range("A66").Value = "mens compression pants pack"
$ws.Range("A67").Value = "boys running pants"
$ws.Range("A68").Value = "knee pads thigh support"
$ws.Range("A69").Value = "youth baseball pants black"
$ws.Range("A70").Value = "knee pads biking adult"
$ws.Range("A71").Value = "youth boys leggings"
$ws.Range("A72").Value = "adult pants"
$ws.Range("A73").Value = "youth baseball compression sleeves"
$ws.Range("A74").Value = "calf sleeves for men football"
$ws.Range("A75").Value = "padded knee sleeve"
$ws.Range("A76").Value = "knee pad exercise"
$ws.Range("A77").Value = "recovery pants men"
$ws.Range("A78").Value = "mens tight"
$ws.Range("A79").Value = "mens outdoor basketball"
$ws.Range("A80").Value = "soccer pants youth"
$ws.Range("A81").Value = "protective shorts"
$ws.Range("A82").Value = "baseball hand guard"
$ws.Range("A83").Value = "bee pants"
$ws.Range("A84").Value = "mens protective pads"
$ws.Range("A85").Value = "mens shorts above knee"
$ws.Range("A86").Value = "basketball knee support for men"
$ws.Range("A87").Value = "compression pads for surgery"
$ws.Range("A88").Value = "snowboarding pants boys"
$ws.Range("A89").Value = "basketball pants for girls"
$ws.Range("A90").Value = "youth baseball pants girls"
$ws.Range("A91").Value = "knee pads for exercise"
$ws.Range("A92").Value = "girls compression knee sleeves"
$ws.Range("A93").Value = "men yoga pant"
$ws.Range("A94").Value = "yoga pants mens black"
$ws.Range("A95").Value = "basketball aids"
$ws.Range("A96").Value = "knee compression sleeve volleyball"
$ws.Range("A97").Value = "impact advanced recovery"
$ws.Range("A98").Value = "best basketball"
$ws.Range("A99").Value = "baseball compression sleeve youth"
$ws.Range("A100").Value = "compression knee pads basketball men"
